# Populate the "column_ID_SI" (column F) values that map each row to its
# corresponding entry/panel in the paper's Supplementary Information, and
# scroll/re-select the sheet view back to the top.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that just need the SI identifier written into column F (no fill).
$ws.Range("F14").Value = "2015_Kieser_GI_1"
$ws.Range("F15").Value = "2016_Nambi"
$ws.Range("F37").Value = "2018_Carey_1A"
$ws.Range("F38").Value = "2018_Carey_1B"
$ws.Range("F39").Value = "2018_Carey_1C"
$ws.Range("F40").Value = "2018_Carey_1D"
$ws.Range("F41").Value = "2018_Carey_1E"
$ws.Range("F42").Value = "2018_Carey_1F"
$ws.Range("F43").Value = "2018_Carey_1G"
$ws.Range("F44").Value = "2018_Carey_1H"
$ws.Range("F45").Value = "2018_Rittershaus_1B"
$ws.Range("F46").Value = "2018_Rittershaus_1A"

# Rows that get the SI identifier AND a yellow highlight fill in column F.
$ws.Range("F27").Value = "2017_Xu_1A"
$ws.Range("F27").Interior.Color = 62207
$ws.Range("F27").Interior.PatternColor = 65535

$ws.Range("F28").Value = "2017_Xu_1B"
$ws.Range("F28").Interior.Color = 62207
$ws.Range("F28").Interior.PatternColor = 65535

$ws.Range("F29").Value = "2017_Xu_1C"
$ws.Range("F29").Interior.Color = 62207
$ws.Range("F29").Interior.PatternColor = 65535

$ws.Range("F31").Value = "2017_Xu_1D"
$ws.Range("F31").Interior.Color = 62207
$ws.Range("F31").Interior.PatternColor = 65535

$ws.Range("F33").Value = "2017_Xu_1E"
$ws.Range("F33").Interior.Color = 62207
$ws.Range("F33").Interior.PatternColor = 65535

# Rows that only get the yellow highlight fill in column F (no text yet).
$ws.Range("F30").Interior.Color = 62207
$ws.Range("F30").Interior.PatternColor = 65535

$ws.Range("F32").Interior.Color = 62207
$ws.Range("F32").Interior.PatternColor = 65535

# Scroll the frozen pane back up to the top and re-select the cell that was
# previously selected relative to the (now-shifted) view.
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 2
$ws.Range("F47").Select()
